$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the attendance date as plain text (e.g. "28-07-2022").
# Some of the new values are day/month-ambiguous (day <= 12), so Excel would
# otherwise auto-parse them into date serials. Force text format first,
# assign the literal string, then restore the default "Normal" style so the
# cell keeps its original (unstyled) appearance while the value stays text.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("A3").Style = "Normal"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("A4").Style = "Normal"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("A5").Style = "Normal"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "08-08-2022"
$ws.Range("A6").Style = "Normal"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "11-08-2022"
$ws.Range("A7").Style = "Normal"

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A8").Style = "Normal"

$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A9").Style = "Normal"

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A10").Style = "Normal"

$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A11").Style = "Normal"

$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A12").Style = "Normal"
$ws.Range("D12").Value = 1
$ws.Range("G12").Value = 1

$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A13").Style = "Normal"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A14").Style = "Normal"

$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A15").Style = "Normal"

$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A16").Style = "Normal"

$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A17").Style = "Normal"

$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A18").Style = "Normal"

$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A19").Style = "Normal"

$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A20").Style = "Normal"

$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "29-09-2022"
$ws.Range("A21").Style = "Normal"

